$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-7) holds a date value that was updated from 2023-09-14
# (serial 45183) to 2023-09-15 (serial 45184). Update each cell's value
# while preserving its existing date formatting/style.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45184
}
